$wb = $excel.ActiveWorkbook

# --- About sheet updates ---
$about = $wb.Worksheets.Item("About")

# Update the "last updated" date in C1 from 1/3/2024 to 3/28/2024
$about.Range("C1").Value = (Get-Date -Year 2024 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0)

# --- FPIEBP sheet updates ---
$fpiebp = $wb.Worksheets.Item("FPIEBP")

# Update priority values for "hard coal" (row 3): production/imports/exports
$fpiebp.Range("B3").Value = 1
$fpiebp.Range("C3").Value = 3
$fpiebp.Range("D3").Value = 2

# Make FPIEBP the active sheet and update its selected cell
$fpiebp.Activate()
$fpiebp.Range("E3").Select()
